# Avances del análisis de texto en Español para la short description
#
# The sheet is a word -> frequency table:
#   column A = rank (0-based), column B = word, column C = frequency count.
# Five Spanish stop/odd words that leaked into the vocabulary ("para",
# "país", "visión", "rayón", "michoacán") are removed from the corpus.
# Removing their rows shifts every following row up by one, which is why
# column A (a literal cached rank, not a formula) has to be recomputed
# afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1
$wordsToRemove = @("para", "país", "visión", "rayón", "michoacán")

# Locate the row for each word with an exact (whole-cell) match so we don't
# accidentally hit substrings like "preparation" or "separation".
$rowsToDelete = @()
foreach ($word in $wordsToRemove) {
    $hit = $ws.Columns("B").Find($word, [Type]::Missing, [Type]::Missing, $xlWhole)
    if ($hit -ne $null) {
        $rowsToDelete += $hit.Row
    }
}

# Delete bottom-to-top so earlier row numbers stay valid while we work.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($rowNum in $rowsToDelete) {
    $ws.Rows($rowNum).Delete()
}

# Recompute the rank column (A2:A<lastRow>) as a plain 0..n-1 sequence now
# that the five rows are gone.
$lastRow = $ws.UsedRange.Rows.Count
$count = $lastRow - 1
$ranks = New-Object 'object[,]' $count,1
for ($i = 0; $i -lt $count; $i++) {
    $ranks[$i, 0] = $i
}
$ws.Range("A2:A$lastRow").Value = $ranks
